{"js": "const replacements = [\n  [\"2023-03-05 Sunday\", \"2023-03-06 Monday\"],\n  [\"34+55=\", \"79-69=\"],\n  [\"12+62=\", \"95-33=\"],\n  [\"16+18=\", \"20+58=\"],\n  [\"30-19=\", \"89-3=\"],\n  [\"50+19=\", \"84-75=\"],\n  [\"73-32=\", \"62+17=\"],\n  [\"56-29=\", \"6+43=\"],\n  [\"21+13=\", \"39-7=\"],\n  [\"49+12=\", \"59-30=\"],\n  [\"43+37=\", \"14+23=\"],\n  [\"4+65=\", \"48-17=\"],\n  [\"17-3=\", \"95-79=\"],\n  [\"89+10=\", \"7+40=\"],\n  [\"11+66=\", \"49-24=\"],\n  [\"27+28=\", \"17+39=\"],\n  [\"64-25=\", \"78-70=\"],\n  [\"20+24=\", \"84-4=\"],\n  [\"71-3=\", \"26+68=\"],\n  [\"36+59=\", \"98-27=\"],\n  [\"42+5=\", \"77-58=\"],\n  [\"83-55=\", \"8+78=\"],\n  [\"66-25=\", \"74-27=\"],\n  [\"4+50=\", \"4+3=\"],\n  [\"2+67=\", \"32+12=\"],\n  [\"3+83=\", \"20+53=\"],\n  [\"65-55=\", \"24+53=\"],\n  [\"39-10=\", \"11+3=\"],\n  [\"25+12=\", \"97+2=\"],\n  [\"38+10=\", \"58+7=\"],\n  [\"15-2=\", \"89-6=\"],\n  [\"98-79=\", \"59-24=\"],\n  [\"6+87=\", \"75-72=\"],\n  [\"0+80=\", \"62-25=\"],\n  [\"94-83=\", \"92-5=\"],\n  [\"38+38=\", \"8+56=\"],\n  [\"80-31=\", \"10+30=\"],\n  [\"66-33=\", \"61-48=\"],\n  [\"36+20=\", \"74-17=\"],\n  [\"56-36=\", \"80-52=\"],\n  [\"22+59=\", \"43-42=\"],\n  [\"20+3=\", \"1+24=\"],\n  [\"29-0=\", \"51-48=\"],\n  [\"25-3=\", \"14+35=\"],\n  [\"71-61=\", \"23-1=\"],\n  [\"8+88=\", \"2+61=\"],\n  [\"21-7=\", \"84-28=\"],\n  [\"21+78=\", \"95-83=\"],\n  [\"64+24=\", \"94-79=\"],\n  [\"85-70=\", \"16-13=\"],\n  [\"89+1=\", \"96-1=\"],\n  [\"49+2=\", \"85+1=\"],\n  [\"62-30=\", \"55+2=\"],\n  [\"16-15=\", \"24+62=\"],\n  [\"80-54=\", \"92-81=\"],\n  [\"25-14=\", \"33+49=\"],\n  [\"76+2=\", \"98-28=\"],\n  [\"7+88=\", \"30+60=\"],\n  [\"92-19=\", \"15+8=\"],\n  [\"18+0=\", \"41-9=\"],\n  [\"54+5=\", \"85-62=\"],\n  [\"85+2=\", \"7+79=\"],\n  [\"52+45=\", \"89-7=\"],\n  [\"54-39=\", \"44+28=\"],\n  [\"22-19=\", \"56+41=\"],\n  [\"51-38=\", \"91-80=\"],\n  [\"80-11=\", \"46-5=\"],\n  [\"21+60=\", \"98-67=\"],\n  [\"37+35=\", \"8-3=\"],\n  [\"27+43=\", \"17+30=\"],\n  [\"1+92=\", \"63-7=\"],\n  [\"64+21=\", \"70-11=\"],\n  [\"5+41=\", \"53+42=\"],\n  [\"55+16=\", \"84-50=\"],\n  [\"10+77=\", \"7+80=\"],\n  [\"7+33=\", \"0+74=\"],\n  [\"87-21=\", \"68-38=\"],\n  [\"25-11=\", \"28+19=\"],\n  [\"20+44=\", \"40+19=\"],\n  [\"83-19=\", \"1+3=\"],\n  [\"17+78=\", \"92-36=\"],\n  [\"90-61=\", \"84-20=\"],\n  [\"78-32=\", \"3+50=\"],\n  [\"52+42=\", \"69-37=\"],\n  [\"19+52=\", \"30-7=\"],\n  [\"67-57=\", \"27+49=\"],\n  [\"55+31=\", \"90-22=\"],\n  [\"9+75=\", \"52-22=\"],\n  [\"22+6=\", \"49+14=\"],\n  [\"12+53=\", \"82-2=\"],\n  [\"90-41=\", \"14+22=\"],\n  [\"58+21=\", \"64-13=\"],\n  [\"5+4=\", \"51-20=\"],\n  [\"95-71=\", \"65-41=\"],\n  [\"60+34=\", \"69-30=\"],\n  [\"23+13=\", \"39+3=\"],\n  [\"35-29=\", \"37-16=\"],\n  [\"16+29=\", \"21+2=\"],\n  [\"2+75=\", \"84-61=\"],\n  [\"55-9=\", \"20+28=\"],\n  [\"17+75=\", \"67-45=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-03-05 Sunday\", \"2023-03-06 Monday\"),\n    @(\"34+55=\", \"79-69=\"),\n    @(\"12+62=\", \"95-33=\"),\n    @(\"16+18=\", \"20+58=\"),\n    @(\"30-19=\", \"89-3=\"),\n    @(\"50+19=\", \"84-75=\"),\n    @(\"73-32=\", \"62+17=\"),\n    @(\"56-29=\", \"6+43=\"),\n    @(\"21+13=\", \"39-7=\"),\n    @(\"49+12=\", \"59-30=\"),\n    @(\"43+37=\", \"14+23=\"),\n    @(\"4+65=\", \"48-17=\"),\n    @(\"17-3=\", \"95-79=\"),\n    @(\"89+10=\", \"7+40=\"),\n    @(\"11+66=\", \"49-24=\"),\n    @(\"27+28=\", \"17+39=\"),\n    @(\"64-25=\", \"78-70=\"),\n    @(\"20+24=\", \"84-4=\"),\n    @(\"71-3=\", \"26+68=\"),\n    @(\"36+59=\", \"98-27=\"),\n    @(\"42+5=\", \"77-58=\"),\n    @(\"83-55=\", \"8+78=\"),\n    @(\"66-25=\", \"74-27=\"),\n    @(\"4+50=\", \"4+3=\"),\n    @(\"2+67=\", \"32+12=\"),\n    @(\"3+83=\", \"20+53=\"),\n    @(\"65-55=\", \"24+53=\"),\n    @(\"39-10=\", \"11+3=\"),\n    @(\"25+12=\", \"97+2=\"),\n    @(\"38+10=\", \"58+7=\"),\n    @(\"15-2=\", \"89-6=\"),\n    @(\"98-79=\", \"59-24=\"),\n    @(\"6+87=\", \"75-72=\"),\n    @(\"0+80=\", \"62-25=\"),\n    @(\"94-83=\", \"92-5=\"),\n    @(\"38+38=\", \"8+56=\"),\n    @(\"80-31=\", \"10+30=\"),\n    @(\"66-33=\", \"61-48=\"),\n    @(\"36+20=\", \"74-17=\"),\n    @(\"56-36=\", \"80-52=\"),\n    @(\"22+59=\", \"43-42=\"),\n    @(\"20+3=\", \"1+24=\"),\n    @(\"29-0=\", \"51-48=\"),\n    @(\"25-3=\", \"14+35=\"),\n    @(\"71-61=\", \"23-1=\"),\n    @(\"8+88=\", \"2+61=\"),\n    @(\"21-7=\", \"84-28=\"),\n    @(\"21+78=\", \"95-83=\"),\n    @(\"64+24=\", \"94-79=\"),\n    @(\"85-70=\", \"16-13=\"),\n    @(\"89+1=\", \"96-1=\"),\n    @(\"49+2=\", \"85+1=\"),\n    @(\"62-30=\", \"55+2=\"),\n    @(\"16-15=\", \"24+62=\"),\n    @(\"80-54=\", \"92-81=\"),\n    @(\"25-14=\", \"33+49=\"),\n    @(\"76+2=\", \"98-28=\"),\n    @(\"7+88=\", \"30+60=\"),\n    @(\"92-19=\", \"15+8=\"),\n    @(\"18+0=\", \"41-9=\"),\n    @(\"54+5=\", \"85-62=\"),\n    @(\"85+2=\", \"7+79=\"),\n    @(\"52+45=\", \"89-7=\"),\n    @(\"54-39=\", \"44+28=\"),\n    @(\"22-19=\", \"56+41=\"),\n    @(\"51-38=\", \"91-80=\"),\n    @(\"80-11=\", \"46-5=\"),\n    @(\"21+60=\", \"98-67=\"),\n    @(\"37+35=\", \"8-3=\"),\n    @(\"27+43=\", \"17+30=\"),\n    @(\"1+92=\", \"63-7=\"),\n    @(\"64+21=\", \"70-11=\"),\n    @(\"5+41=\", \"53+42=\"),\n    @(\"55+16=\", \"84-50=\"),\n    @(\"10+77=\", \"7+80=\"),\n    @(\"7+33=\", \"0+74=\"),\n    @(\"87-21=\", \"68-38=\"),\n    @(\"25-11=\", \"28+19=\"),\n    @(\"20+44=\", \"40+19=\"),\n    @(\"83-19=\", \"1+3=\"),\n    @(\"17+78=\", \"92-36=\"),\n    @(\"90-61=\", \"84-20=\"),\n    @(\"78-32=\", \"3+50=\"),\n    @(\"52+42=\", \"69-37=\"),\n    @(\"19+52=\", \"30-7=\"),\n    @(\"67-57=\", \"27+49=\"),\n    @(\"55+31=\", \"90-22=\"),\n    @(\"9+75=\", \"52-22=\"),\n    @(\"22+6=\", \"49+14=\"),\n    @(\"12+53=\", \"82-2=\"),\n    @(\"90-41=\", \"14+22=\"),\n    @(\"58+21=\", \"64-13=\"),\n    @(\"5+4=\", \"51-20=\"),\n    @(\"95-71=\", \"65-41=\"),\n    @(\"60+34=\", \"69-30=\"),\n    @(\"23+13=\", \"39+3=\"),\n    @(\"35-29=\", \"37-16=\"),\n    @(\"16+29=\", \"21+2=\"),\n    @(\"2+75=\", \"84-61=\"),\n    @(\"55-9=\", \"20+28=\"),\n    @(\"17+75=\", \"67-45=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
